$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.425.61'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.92%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.801.52'
$ws.Range("D3").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.68%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '608.51'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.81%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '163.65'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.46%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.799.30'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.10%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("E9").Value = '  -0.41%  '

$ws.Range("E10").Value = '  -0.22%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.98'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +10.45%  '

$ws.Range("E12").Value = '  -0.49%  '

$ws.Range("E13").Value = '  -1.47%  '

$ws.Range("E14").Value = '  -2.51%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.439.72'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.16%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.806.07'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.39%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.452.92'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.91%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.06'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.07%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.114'
$ws.Range("D19").Style = "Normal"

$ws.Range("E20").Value = '  -0.05%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '461.91'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.46%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.59'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.94%  '

$ws.Range("E23").Value = '  -0.57%  '

$ws.Range("E24").Value = '  -0.24%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.47'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.20%  '

$ws.Range("E26").Value = '  -1.12%  '

$ws.Range("E27").Value = '  -0.86%  '

$ws.Range("E28").Value = '  +0.03%  '

$ws.Range("E29").Value = '  -0.51%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.946.61'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.26%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.63'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.84%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.20'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.16%  '

$ws.Range("E33").Value = '  -2.15%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.04'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.54%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.04%  '

$ws.Range("E36").Value = '  -0.38%  '

$ws.Range("E37").Value = '  +0.47%  '

$ws.Range("E38").Value = '  +6.34%  '

$ws.Range("E39").Value = '  +1.00%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.979'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.65%  '

$ws.Range("E41").Value = '  -0.96%  '

$ws.Range("E42").Value = '  -0.03%  '

$ws.Range("B44").Value = 'TheGraph'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.297'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.12%  '

$ws.Range("B45").Value = 'Monero'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '153.14'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.58%  '

$ws.Range("B46").Value = 'ONDO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.41'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.91%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '43.03'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.84%  '

$ws.Range("B48").Value = 'OKB'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '46.90'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.71%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.37'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.17%  '

$ws.Range("E50").Value = '  +0.20%  '

$ws.Range("B51").Value = 'Bittensor'
$ws.Range("C51").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '379.36'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.00%  '
